$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.444654226303101
$ws.Range("B1").Value = 2.256736040115356
$ws.Range("C1").Value = 5.070332050323486
$ws.Range("D1").Value = 3.238076686859131
$ws.Range("E1").Value = 1.131497025489807
